# Remove the "Certifications" section's content:
#   - the "Certifications" Heading2 paragraph becomes empty (and gains
#     an explicit left indent of 0)
#   - the "Oracle Database 12c: Basic SQL" body paragraph becomes empty
#   - the blank spacer paragraph right before the heading is deleted
#     entirely (paragraph mark and all)

$d = $word.ActiveDocument

function Find-ParaIndexByText($doc, $pattern) {
    $idx = 0
    foreach ($p in $doc.Paragraphs) {
        $idx = $idx + 1
        if ($p.Range.Text -match $pattern) {
            return $idx
        }
    }
    return -1
}

# Locate the two text-bearing paragraphs by their content rather than a
# hard-coded index, so the script is resilient to unrelated edits earlier
# in the document.
$certIndex = Find-ParaIndexByText $d "Certifications"
$oracleIndex = Find-ParaIndexByText $d "Oracle Database"

# Clear the runs in the "Oracle Database 12c: Basic SQL" paragraph while
# keeping the paragraph (and its paragraph mark / formatting) intact.
# Assigning "" to the text up to (but excluding) the paragraph mark removes
# the run(s) cleanly instead of leaving a stray empty <w:r>.
$oraclePara = $d.Paragraphs.Item($oracleIndex)
$oracleRange = $oraclePara.Range
$d.Range($oracleRange.Start, $oracleRange.End - 1).Text = ""

# Clear the runs in the "Certifications" heading paragraph the same way.
$certPara = $d.Paragraphs.Item($certIndex)
$certRange = $certPara.Range
$d.Range($certRange.Start, $certRange.End - 1).Text = ""

# Give the now-empty heading paragraph an explicit left indent of 0.
$certPara2 = $d.Paragraphs.Item($certIndex)
$certPara2.Range.ParagraphFormat.LeftIndent = 0

# Delete the blank spacer paragraph immediately preceding the heading
# (its own paragraph mark included, so it disappears entirely).
$spacerIndex = $certIndex - 1
$spacerPara = $d.Paragraphs.Item($spacerIndex)
$spacerPara.Range.Delete()
